$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 29: Dripping with Venom / Weak Blinding Potion
$ws.Cells.Item(29, 8).Value2 = 683.5
$ws.Cells.Item(29, 9).Value2 = 773.25
$ws.Cells.Item(29, 10).Value2 = 504
$ws.Cells.Item(29, 11).Value2 = 2319.75
$ws.Cells.Item(29, 12).Value2 = 1512
$ws.Cells.Item(29, 13).Value2 = -2038.75
$ws.Cells.Item(29, 14).Value2 = -2074

# ALC row 32: Automata for the People / Crab Oil
$ws.Cells.Item(32, 8).Value2 = 2528.6667
$ws.Cells.Item(32, 10).Value2 = 1305
$ws.Cells.Item(32, 12).Value2 = 1305
$ws.Cells.Item(32, 14).Value2 = -1957

# ALC row 38: Just Give Him a Serum / Hi-Potion of Strength
$ws.Cells.Item(38, 8).Value2 = 203.75
$ws.Cells.Item(38, 9).Value2 = 102
$ws.Cells.Item(38, 10).Value2 = 509
$ws.Cells.Item(38, 11).Value2 = 306
$ws.Cells.Item(38, 12).Value2 = 1527
$ws.Cells.Item(38, 13).Value2 = 66
$ws.Cells.Item(38, 14).Value2 = -2271

# ALC row 40: Stuck in the Moment / Horn Glue
$ws.Cells.Item(40, 8).Value2 = 1011.75
$ws.Cells.Item(40, 9).Value2 = 978.6667
$ws.Cells.Item(40, 10).Value2 = 1111
$ws.Cells.Item(40, 11).Value2 = 978.6667
$ws.Cells.Item(40, 12).Value2 = 1111
$ws.Cells.Item(40, 13).Value2 = -803.6667
$ws.Cells.Item(40, 14).Value2 = -1461

# ALC row 58: A Matter of Vital Importance / Mega-Potion of Vitality
$ws.Cells.Item(58, 8).Value2 = 2832.3076
$ws.Cells.Item(58, 9).Value2 = 204
$ws.Cells.Item(58, 10).Value2 = 4475
$ws.Cells.Item(58, 11).Value2 = 612
$ws.Cells.Item(58, 12).Value2 = 13425
$ws.Cells.Item(58, 13).Value2 = -462
$ws.Cells.Item(58, 14).Value2 = -13725

# ALC row 87: There Was a Late Fee / Noble Gold
$ws.Cells.Item(87, 8).Value2 = 50000
$ws.Cells.Item(87, 10).Value2 = 50000
$ws.Cells.Item(87, 12).Value2 = 50000
$ws.Cells.Item(87, 14).Value2 = -52496

# ALC row 90: A Gate Arcane Is Dragon's Bane (L) / Noble Gold
$ws.Cells.Item(90, 8).Value2 = 50000
$ws.Cells.Item(90, 10).Value2 = 50000
$ws.Cells.Item(90, 12).Value2 = 150000
$ws.Cells.Item(90, 14).Value2 = -162480

# ALC row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Cells.Item(112, 8).Value2 = 66667944
$ws.Cells.Item(112, 9).Value2 = 646.6667
$ws.Cells.Item(112, 10).Value2 = 83334770
$ws.Cells.Item(112, 11).Value2 = 1940.0001
$ws.Cells.Item(112, 12).Value2 = 250004310
$ws.Cells.Item(112, 13).Value2 = -832.0001
$ws.Cells.Item(112, 14).Value2 = -250006526

# ALC row 129: Practical Command / Commanding Craftsman's Draught
$ws.Cells.Item(129, 8).Value2 = 1381.25
$ws.Cells.Item(129, 9).Value2 = 583.44446
$ws.Cells.Item(129, 10).Value2 = 1647.1852
$ws.Cells.Item(129, 11).Value2 = 1750.33338
$ws.Cells.Item(129, 12).Value2 = 4941.5556
$ws.Cells.Item(129, 13).Value2 = 3249.66662
$ws.Cells.Item(129, 14).Value2 = -14941.5556

# ALC row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Cells.Item(132, 8).Value2 = 266249.7
$ws.Cells.Item(132, 9).Value2 = 297473.06
$ws.Cells.Item(132, 10).Value2 = 851
$ws.Cells.Item(132, 11).Value2 = 892419.1799999999
$ws.Cells.Item(132, 12).Value2 = 2553
$ws.Cells.Item(132, 13).Value2 = -889889.1799999999
$ws.Cells.Item(132, 14).Value2 = -7613

# ALC row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Cells.Item(138, 8).Value2 = 2145.775
$ws.Cells.Item(138, 9).Value2 = 1594.85
$ws.Cells.Item(138, 10).Value2 = 2696.7
$ws.Cells.Item(138, 11).Value2 = 4784.549999999999
$ws.Cells.Item(138, 12).Value2 = 8090.099999999999
$ws.Cells.Item(138, 13).Value2 = 355.4500000000007
$ws.Cells.Item(138, 14).Value2 = -18370.1

$ws = $wb.Worksheets.Item("ARM")
# ARM row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Cells.Item(61, 8).Value2 = 1931.5555
$ws.Cells.Item(61, 9).Value2 = 1171.8334
$ws.Cells.Item(61, 10).Value2 = 8009.3335
$ws.Cells.Item(61, 11).Value2 = 1171.8334
$ws.Cells.Item(61, 12).Value2 = 8009.3335
$ws.Cells.Item(61, 13).Value2 = -959.8334
$ws.Cells.Item(61, 14).Value2 = -8433.333500000001

# ARM row 74: As the Bolt Flies / Titanium Nugget
$ws.Cells.Item(74, 8).Value2 = 3909.476
$ws.Cells.Item(74, 9).Value2 = 883.5
$ws.Cells.Item(74, 10).Value2 = 5422.4644
$ws.Cells.Item(74, 11).Value2 = 883.5
$ws.Cells.Item(74, 12).Value2 = 5422.4644
$ws.Cells.Item(74, 13).Value2 = -9.5
$ws.Cells.Item(74, 14).Value2 = -7170.4644

# ARM row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Cells.Item(77, 8).Value2 = 3909.476
$ws.Cells.Item(77, 9).Value2 = 883.5
$ws.Cells.Item(77, 10).Value2 = 5422.4644
$ws.Cells.Item(77, 11).Value2 = 4417.5
$ws.Cells.Item(77, 12).Value2 = 27112.322
$ws.Cells.Item(77, 13).Value2 = -49.5
$ws.Cells.Item(77, 14).Value2 = -35848.322

# ARM row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Cells.Item(132, 8).Value2 = 1382370.5
$ws.Cells.Item(132, 9).Value2 = 1726913.2
$ws.Cells.Item(132, 10).Value2 = 4199.6
$ws.Cells.Item(132, 11).Value2 = 5180739.6
$ws.Cells.Item(132, 12).Value2 = 12598.8
$ws.Cells.Item(132, 13).Value2 = -5178209.6
$ws.Cells.Item(132, 14).Value2 = -17658.8

# ARM row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Cells.Item(136, 8).Value2 = 1931.5555
$ws.Cells.Item(136, 9).Value2 = 1171.8334
$ws.Cells.Item(136, 10).Value2 = 8009.3335
$ws.Cells.Item(136, 11).Value2 = 3515.5002
$ws.Cells.Item(136, 12).Value2 = 24028.0005
$ws.Cells.Item(136, 13).Value2 = -965.5001999999999
$ws.Cells.Item(136, 14).Value2 = -29128.0005

$ws = $wb.Worksheets.Item("BSM")
# BSM row 64: With Bearings Straight / Mythrite Nugget
$ws.Cells.Item(64, 8).Value2 = 1482.72
$ws.Cells.Item(64, 10).Value2 = 762.2222
$ws.Cells.Item(64, 12).Value2 = 762.2222
$ws.Cells.Item(64, 14).Value2 = -1212.2222

# BSM row 67: Bearing the Brunt (L) / Mythrite Nugget
$ws.Cells.Item(67, 8).Value2 = 1482.72
$ws.Cells.Item(67, 10).Value2 = 762.2222
$ws.Cells.Item(67, 12).Value2 = 762.2222
$ws.Cells.Item(67, 14).Value2 = -2322.2222

# BSM row 107: The Gold Experience / Deepgold Nugget
$ws.Cells.Item(107, 8).Value2 = 224828.19
$ws.Cells.Item(107, 9).Value2 = 314183.53
$ws.Cells.Item(107, 10).Value2 = 1439.8334
$ws.Cells.Item(107, 11).Value2 = 314183.53
$ws.Cells.Item(107, 12).Value2 = 1439.8334
$ws.Cells.Item(107, 13).Value2 = -312263.53
$ws.Cells.Item(107, 14).Value2 = -5279.8334

# BSM row 132: Always Be Prepaired / Mountain Chromite Twinfangs
$ws.Cells.Item(132, 8).Value2 = 42000
$ws.Cells.Item(132, 10).Value2 = 42000
$ws.Cells.Item(132, 12).Value2 = 42000
$ws.Cells.Item(132, 14).Value2 = -52120

# BSM row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Cells.Item(134, 8).Value2 = 41478.867
$ws.Cells.Item(134, 9).Value2 = 47570.69
$ws.Cells.Item(134, 11).Value2 = 142712.07
$ws.Cells.Item(134, 13).Value2 = -140177.07

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31: Wall Not Found / Walnut Lumber
$ws.Cells.Item(31, 8).Value2 = 1666.1538
$ws.Cells.Item(31, 10).Value2 = 2800
$ws.Cells.Item(31, 12).Value2 = 2800
$ws.Cells.Item(31, 14).Value2 = -3390

# CRP row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Cells.Item(34, 8).Value2 = 1666.1538
$ws.Cells.Item(34, 10).Value2 = 2800
$ws.Cells.Item(34, 12).Value2 = 2800
$ws.Cells.Item(34, 14).Value2 = -3204

# CRP row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Cells.Item(134, 8).Value2 = 2112.1904
$ws.Cells.Item(134, 9).Value2 = 2233.3057
$ws.Cells.Item(134, 10).Value2 = 1385.5
$ws.Cells.Item(134, 11).Value2 = 6699.9171
$ws.Cells.Item(134, 12).Value2 = 4156.5
$ws.Cells.Item(134, 13).Value2 = -4164.9171
$ws.Cells.Item(134, 14).Value2 = -9226.5

$ws = $wb.Worksheets.Item("CUL")
# CUL row 107: Slippery Service / Frantoio Oil
$ws.Cells.Item(107, 8).Value2 = 57476.6
$ws.Cells.Item(107, 9).Value2 = 62657.812
$ws.Cells.Item(107, 10).Value2 = 53113.473
$ws.Cells.Item(107, 11).Value2 = 187973.436
$ws.Cells.Item(107, 12).Value2 = 159340.419
$ws.Cells.Item(107, 13).Value2 = -186053.436
$ws.Cells.Item(107, 14).Value2 = -163180.419

# CUL row 113: Can't Eat Just One / Night Vinegar
$ws.Cells.Item(113, 8).Value2 = 533.9783
$ws.Cells.Item(113, 10).Value2 = 539.71875
$ws.Cells.Item(113, 12).Value2 = 1619.15625
$ws.Cells.Item(113, 14).Value2 = -5959.15625

$ws = $wb.Worksheets.Item("LTW")
# LTW row 93: Hide to Go Seek / Gagana Leather
$ws.Cells.Item(93, 8).Value2 = 2390.9092
$ws.Cells.Item(93, 9).Value2 = 2140
$ws.Cells.Item(93, 10).Value2 = 2600
$ws.Cells.Item(93, 11).Value2 = 2140
$ws.Cells.Item(93, 12).Value2 = 2600
$ws.Cells.Item(93, 13).Value2 = -892
$ws.Cells.Item(93, 14).Value2 = -5096

# LTW row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Cells.Item(132, 8).Value2 = 104327.37
$ws.Cells.Item(132, 9).Value2 = 126178.11
$ws.Cells.Item(132, 10).Value2 = 5999
$ws.Cells.Item(132, 11).Value2 = 378534.33
$ws.Cells.Item(132, 12).Value2 = 17997
$ws.Cells.Item(132, 13).Value2 = -376004.33
$ws.Cells.Item(132, 14).Value2 = -23057
